$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new timesheet rows (83, 84, 85)
$ws.Range("A83").Value2 = 42921
$ws.Range("B83").Value2 = "Revize"
$ws.Range("C83").Value2 = 5

$ws.Range("A84").Value2 = 42922
$ws.Range("C84").Value2 = 8

$ws.Range("A85").Value2 = 42923
$ws.Range("C85").Value2 = 3

# Update the selection to match the new last row, as in the diff
$ws.Range("A85:C85").Select()
